$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7 (shifts rows 7-18 down to 8-19)
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the "4-Wire Cable" entry.
# Write order matches how the new shared strings ended up ordered in the
# workbook: link URL, then notes, then the part name.
$ws.Cells.Item(7, 5).Value = "https://www.amazon.com/gp/product/B08JTZCJV1/ref=ppx_yo_dt_b_search_asin_title?ie=UTF8&psc=1"
$ws.Cells.Item(7, 5).Style = "Hyperlink"
$ws.Cells.Item(7, 6).Value = "Need appx. 1 meter"
$ws.Cells.Item(7, 2).Value = "4-Wire Cable"

$ws.Cells.Item(7, 1).Value = 3
$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(7, 4).Value = "Amazon"

# Renumber the "Number" column so it stays sequential after the insert
for ($r = 8; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 4
}

$ws.Range("F7").Select()
